$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45179 = 2023-09-10) for every
# data row (2-171). Update it to 45180 (2023-09-11) to reflect the new "changed" date.
for ($r = 2; $r -le 171; $r++) {
    $ws.Cells.Item($r, 3).Value = 45180
}
